$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Translate header row labels from English to French
$ws.Range("A1").Value = "auteur"
$ws.Range("B1").Value = "codepostal"
$ws.Range("C1").Value = "ville"
$ws.Range("D1").Value = "élément"
$ws.Range("E1").Value = "élément2"
$ws.Range("G1").Value = "titre"

# Move the active cell / selection to K4
$ws.Range("K4").Select()
